# "TARUGO ESPIGA DISMAY" price list (Hoja1) refresh:
#   - bump the price-list date (A1) forward one month
#   - update the four bag prices in column D (rows 32-35)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436

$ws.Range("D32").Value = 13166.095
$ws.Range("D33").Value = 10457.641
$ws.Range("D34").Value = 10432.562
$ws.Range("D35").Value = 15009.349
